# Update gh-pages output data (想去人数 column) regenerated at 456a3b4.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 292
$ws1.Range("F4").Value = 16719
$ws1.Range("F5").Value = 27
$ws1.Range("F6").Value = 1635
$ws1.Range("F12").Value = 11607
$ws1.Range("F14").Value = 1282
$ws1.Range("F15").Value = 4594
$ws1.Range("F18").Value = 65
$ws1.Range("F19").Value = 885

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 292
$ws4.Range("F5").Value = 16719
$ws4.Range("F6").Value = 27
$ws4.Range("F7").Value = 1635
$ws4.Range("F15").Value = 11607
$ws4.Range("F17").Value = 1282
$ws4.Range("F18").Value = 4594
$ws4.Range("F21").Value = 65
$ws4.Range("F22").Value = 885
